# Apply the crypto price/volume refresh described by the commit diff.
# Each D/E cell is rewritten with its new scraped text; B41:E42 additionally
# swap the PaxDollar/mCoin rows. Values that look like plain numbers (a single
# decimal point, e.g. "215.36") are written with a leading quote so Excel keeps
# storing them as text (matching the original inline-string cells) instead of
# silently converting them to numbers and dropping significant trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.833.17"
$ws.Range("E2").Value = "  -1.33%  "

# Row 3
$ws.Range("D3").Value = "1.636.12"
$ws.Range("E3").Value = "  -1.28%  "

# Row 4
$ws.Range("E4").Value = "  -0.30%  "

# Row 5
$ws.Range("D5").Value = "`'215.36"
$ws.Range("E5").Value = "  -0.41%  "

# Row 6
$ws.Range("D6").Value = "`'0.5024"
$ws.Range("E6").Value = "  -2.14%  "

# Row 8
$ws.Range("D8").Value = "`'0.2572"
$ws.Range("E8").Value = "  -0.45%  "

# Row 9
$ws.Range("D9").Value = "`'0.06410"
$ws.Range("E9").Value = "  -0.23%  "

# Row 10
$ws.Range("D10").Value = "`'19.64"
$ws.Range("E10").Value = "  -1.77%  "

# Row 11
$ws.Range("D11").Value = "`'0.07708"
$ws.Range("E11").Value = "  -0.98%  "

# Row 12
$ws.Range("D12").Value = "`'4.245"

# Row 13
$ws.Range("D13").Value = "1.637.07"
$ws.Range("E13").Value = "  -1.41%  "

# Row 14
$ws.Range("D14").Value = "1.861.34"
$ws.Range("E14").Value = "  -1.26%  "

# Row 15
$ws.Range("D15").Value = "`'0.5439"
$ws.Range("E15").Value = "  -1.81%  "

# Row 16
$ws.Range("D16").Value = "0.0₅7934"
$ws.Range("E16").Value = "  -1.37%  "

# Row 17
$ws.Range("D17").Value = "`'63.49"
$ws.Range("E17").Value = "  -1.17%  "

# Row 18
$ws.Range("D18").Value = "25.867.41"
$ws.Range("E18").Value = "  -1.32%  "

# Row 19
$ws.Range("E19").Value = "  -0.16%  "

# Row 20
$ws.Range("D20").Value = "`'203.00"
$ws.Range("E20").Value = "  -3.71%  "

# Row 21
$ws.Range("D21").Value = "`'4.328"
$ws.Range("E21").Value = "  -1.98%  "

# Row 22
$ws.Range("D22").Value = "`'9.949"
$ws.Range("E22").Value = "  -1.18%  "

# Row 23
$ws.Range("D23").Value = "`'5.980"
$ws.Range("E23").Value = "  -0.52%  "

# Row 25
$ws.Range("D25").Value = "`'1.922"
$ws.Range("E25").Value = "  +11.47%  "

# Row 26
$ws.Range("D26").Value = "`'141.17"
$ws.Range("E26").Value = "  -2.06%  "

# Row 27
$ws.Range("D27").Value = "`'0.1144"
$ws.Range("E27").Value = "  -2.25%  "

# Row 28
$ws.Range("E28").Value = "  -0.43%  "

# Row 29
$ws.Range("E29").Value = "  -3.78%  "

# Row 30
$ws.Range("D30").Value = "`'1.240"
$ws.Range("E30").Value = "  -0.68%  "

# Row 31
$ws.Range("D31").Value = "`'0.05010"
$ws.Range("E31").Value = "  -2.23%  "

# Row 32
$ws.Range("D32").Value = "`'3.264"
$ws.Range("E32").Value = "  -2.59%  "

# Row 33
$ws.Range("E33").Value = "  -1.63%  "

# Row 34
$ws.Range("D34").Value = "`'1.536"
$ws.Range("E34").Value = "  -1.97%  "

# Row 35
$ws.Range("D35").Value = "`'2.364"
$ws.Range("E35").Value = "  -0.34%  "

# Row 36
$ws.Range("D36").Value = "1.175.94"
$ws.Range("E36").Value = "  +1.12%  "

# Row 37
$ws.Range("D37").Value = "`'0.8939"
$ws.Range("E37").Value = "  -3.87%  "

# Row 38
$ws.Range("D38").Value = "`'2.606"
$ws.Range("E38").Value = "  -5.39%  "

# Row 39
$ws.Range("D39").Value = "`'0.5614"
$ws.Range("E39").Value = "  -1.56%  "

# Row 40
$ws.Range("D40").Value = "`'0.01559"
$ws.Range("E40").Value = "  -2.11%  "

# Row 41
$ws.Range("B41").Value = "mCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D41").Value = "`'2.554"
$ws.Range("E41").Value = "  -0.45%  "

# Row 42
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "`'1.002"
$ws.Range("E42").Value = "  -0.25%  "

# Row 43
$ws.Range("D43").Value = "`'5.681"
$ws.Range("E43").Value = "  +0.36%  "

# Row 44
$ws.Range("D44").Value = "`'0.8074"
$ws.Range("E44").Value = "  -3.79%  "

# Row 45
$ws.Range("D45").Value = "`'99.35"
$ws.Range("E45").Value = "  -1.04%  "

# Row 46
$ws.Range("D46").Value = "1.773.18"
$ws.Range("E46").Value = "  -1.23%  "

# Row 47
$ws.Range("E47").Value = "  -0.27%  "

# Row 48
$ws.Range("D48").Value = "`'0.4516"
$ws.Range("E48").Value = "  -0.58%  "

# Row 49
$ws.Range("D49").Value = "`'1.004"
$ws.Range("E49").Value = "  +0.10%  "

# Row 50
$ws.Range("D50").Value = "`'54.76"
$ws.Range("E50").Value = "  -1.93%  "

# Row 51
$ws.Range("D51").Value = "`'0.05063"
$ws.Range("E51").Value = "  +0.10%  "
